$d = $word.ActiveDocument
$content = $d.Content

# The paragraph currently contains three runs:
#   1. "<id>"      (Courier New, color 7f6000, sz/szCs 18)
#   2. "p057v_2"   (color 000000)
#   3. "</id>"     (Courier New, color 7f6000, sz/szCs 18)
# They need to be merged into a single run "<id>p057v_2</id>" that keeps
# the formatting of the first run. We do this by deleting the text of the
# 2nd+3rd runs and then re-typing it right after the (untouched) first run,
# so Word extends/merges it into that run instead of creating new runs.

$tail = $content.Duplicate
$tail.Find.ClearFormatting()
$tail.Find.Execute("p057v_2</id>", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $tail.Find.Found) {
    throw "Could not find the text 'p057v_2</id>' to merge."
}
$tail.Delete()

$head = $content.Duplicate
$head.Find.ClearFormatting()
$head.Find.Execute("<id>", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $head.Find.Found) {
    throw "Could not find the '<id>' run to append to."
}
$head.Collapse(0)
$head.InsertAfter("p057v_2</id>")
